$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: fill in "Gedaan" (A) and "Duur" (I) text; date (J18) already present
$ws.Range("A18").Value = "Vraagpagina"
$ws.Range("I18").Value = "3 uur"

# Row 19: Groesgesprek met Robert / 0,5 uur / 19-3-2012
$ws.Range("A19").Value = "Groesgesprek met Robert"
$ws.Range("I19").Value = "0,5 uur"
$ws.Range("J19").Value = 40987

# Row 20: Vraagpagina afgemaakt / 2 uur / 19-3-2012
$ws.Range("A20").Value = "Vraagpagina afgemaakt"
$ws.Range("I20").Value = "2 uur"
$ws.Range("J20").Value = 40987

# Row 21: Tagpagina  (trailing space) / 2 uur / 20-3-2012
$ws.Range("A21").Value = "Tagpagina "
$ws.Range("I21").Value = "2 uur"
$ws.Range("J21").Value = 40988

# Row 22: Tagpagina  (trailing space) / afgemaakt / 0,5 uur / 21-3-2012
$ws.Range("A22").Value = "Tagpagina "
$ws.Range("B22").Value = "afgemaakt"
$ws.Range("I22").Value = "0,5 uur"
$ws.Range("J22").Value = 40989

# Copy the date format (numFmtId 14, m/d/yyyy style) from J18 onto the newly-dated cells
$ws.Range("J18").Copy()
$ws.Range("J19:J22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the sheet view: scroll so row 7 is at top, select B26
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("B26").Select()
